$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.795.95'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '2.549.09'
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''308.83'
$ws.Range("E5").Value = '  -2.98%  '
$ws.Range("D6").Value = '''101.54'
$ws.Range("E6").Value = '  +4.82%  '
$ws.Range("E7").Value = '  -0.45%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").Value = '''36.25'
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("D12").Value = '''7.41'
$ws.Range("E12").Value = '  -1.97%  '
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").Value = '2.940.00'
$ws.Range("E14").Value = '  +0.62%  '
$ws.Range("D15").Value = '''15.96'
$ws.Range("E15").Value = '  +5.46%  '
$ws.Range("D16").Value = '2.548.26'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("E17").Value = '  -1.29%  '
$ws.Range("D18").Value = '42.829.53'
$ws.Range("D19").Value = '''6.77'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").Value = '''12.37'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("E21").Value = '  -1.32%  '
$ws.Range("D22").Value = '''69.43'
$ws.Range("E22").Value = '  -0.39%  '
$ws.Range("D23").Value = '''248.50'
$ws.Range("E23").Value = '  -1.99%  '
$ws.Range("E24").Value = '  -2.49%  '
$ws.Range("E25").Value = '  +0.34%  '
$ws.Range("D26").Value = '''26.55'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").Value = '''40.69'
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("E29").Value = '  -1.91%  '
$ws.Range("E30").Value = '  -4.01%  '
$ws.Range("D31").Value = '''156.57'
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("E32").Value = '  -2.79%  '
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("D34").Value = '''3.30'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("E35").Value = '  -3.00%  '
$ws.Range("E36").Value = '  -2.64%  '
$ws.Range("D37").Value = '''2.62'
$ws.Range("E37").Value = '  +6.24%  '
$ws.Range("D38").Value = '''18.35'
$ws.Range("E38").Value = '  -5.36%  '
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '''4.25'
$ws.Range("E41").Value = '  +11.55%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '''22.57'
$ws.Range("E42").Value = '  +3.76%  '
$ws.Range("E43").Value = '  -0.18%  '
$ws.Range("E44").Value = '  -1.54%  '
$ws.Range("D45").Value = '''3.28'
$ws.Range("E45").Value = '  -0.15%  '
$ws.Range("D46").Value = '1.985.98'
$ws.Range("E46").Value = '  -1.20%  '
$ws.Range("D48").Value = '2.794.12'
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("D49").Value = '''81.45'
$ws.Range("E49").Value = '  -3.33%  '
$ws.Range("E50").Value = '  +0.64%  '
$ws.Range("D51").Value = '''73.67'
$ws.Range("E51").Value = '  -1.66%  '
